$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text updates (rich-text shared strings; all runs share identical rPr) ---
$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"

# --- Crime complaint table updates (rows 14-33) ---
# Row 14
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = "#,##0.0;""-""#,##0.0"
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 120
# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 81.818181818181
$ws.Range("I16").Value = 91
$ws.Range("J16").Value = 69
$ws.Range("K16").Value = 31.884057971014
$ws.Range("L16").Value = 16.666666666666
$ws.Range("M16").Value = 16.666666666666
$ws.Range("N16").Value = -81.466395112016
# Row 17
$ws.Range("F17").Value = 27
$ws.Range("H17").Value = 92.857142857142
$ws.Range("I17").Value = 133
$ws.Range("J17").Value = 99
$ws.Range("K17").Value = 34.343434343434
$ws.Range("L17").Value = 35.714285714285
$ws.Range("M17").Value = 129.310344827586
$ws.Range("N17").Value = 6.4
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = 38.888888888888
$ws.Range("I18").Value = 169
$ws.Range("J18").Value = 91
$ws.Range("K18").Value = 85.714285714285
$ws.Range("L18").Value = 49.557522123893
$ws.Range("M18").Value = 23.357664233576
$ws.Range("N18").Value = -78.221649484536
# Row 19
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -10.526315789473
$ws.Range("F19").Value = 75
$ws.Range("G19").Value = 66
$ws.Range("H19").Value = 13.636363636363
$ws.Range("I19").Value = 509
$ws.Range("J19").Value = 455
$ws.Range("K19").Value = 11.868131868131
$ws.Range("L19").Value = 6.041666666666
$ws.Range("M19").Value = -24.592592592592
$ws.Range("N19").Value = -58.918482647296
# Row 20
$ws.Range("C20").Value = 1
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 1
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("E20").Value = 0
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I20").Value = 23
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 43.75
$ws.Range("L20").Value = -28.125
$ws.Range("M20").Value = 43.75
$ws.Range("N20").Value = -96.068376068376
# Row 21
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = 3.030303030303
$ws.Range("F21").Value = 155
$ws.Range("G21").Value = 113
$ws.Range("H21").Value = 37.16814159292
$ws.Range("I21").Value = 936
$ws.Range("J21").Value = 736
$ws.Range("K21").Value = 27.173913043478
$ws.Range("L21").Value = 15.841584158415
$ws.Range("M21").Value = -3.405572755417
$ws.Range("N21").Value = -71.003717472119
# Row 22
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 34
$ws.Range("J22").Value = 29
$ws.Range("K22").Value = 17.241379310344
$ws.Range("L22").Value = -15
$ws.Range("M22").Value = -5.555555555555
# Row 24
$ws.Range("C24").Value = 48
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = -5.882352941176
$ws.Range("F24").Value = 263
$ws.Range("G24").Value = 209
$ws.Range("H24").Value = 25.837320574162
$ws.Range("I24").Value = 1312
$ws.Range("J24").Value = 1403
$ws.Range("K24").Value = -6.486101211689
$ws.Range("L24").Value = 31.859296482412
$ws.Range("M24").Value = 56.37663885578
# Row 25
$ws.Range("C25").Value = 36
$ws.Range("D25").Value = 49
$ws.Range("E25").Value = -26.530612244898
$ws.Range("F25").Value = 199
$ws.Range("G25").Value = 179
$ws.Range("H25").Value = 11.173184357541
$ws.Range("I25").Value = 1051
$ws.Range("J25").Value = 1200
$ws.Range("K25").Value = -12.416666666666
$ws.Range("L25").Value = 43.775649794801
# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = -44.444444444444
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 61
$ws.Range("H26").Value = -37.704918032786
$ws.Range("I26").Value = 268
$ws.Range("J26").Value = 286
$ws.Range("K26").Value = -6.293706293706
$ws.Range("L26").Value = 17.543859649122
$ws.Range("M26").Value = 35.353535353535
# Row 27
$ws.Range("D27").Value = 1
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 66.666666666666
$ws.Range("L27").Value = 15.384615384615
# Row 28
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 40
$ws.Range("I28").Value = 58
$ws.Range("J28").Value = 52
$ws.Range("K28").Value = 11.538461538461
$ws.Range("L28").Value = -3.333333333333
# Row 33
$ws.Range("C33").Value = 1
$ws.Range("C33").NumberFormat = "#,##0"
$ws.Range("F33").Value = 1
$ws.Range("F33").NumberFormat = "#,##0"
$ws.Range("I33").Value = 1
$ws.Range("I33").NumberFormat = "#,##0"
$ws.Range("L33").Value = -66.666666666666
